$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Reference Lists")

# The first three reference pages (Divisions, Data Set Status, Data Set Type)
# have completed testing - replace their "X" marker with today's completion
# date, and clear the "X" marker for the remaining untested pages.
$ws1.Range("B3:B16").ClearContents()

# Add the two new reference-list rows covering the data set preset fields.
$ws1.Range("A17").Value = "Data Set Status"
$ws1.Range("A18").Value = "Data Set Type"

# Record the completion date (12/4/2023) for the rows that finished testing.
$ws1.Range("B2").NumberFormat = "mm-dd-yy"
$ws1.Range("B2").Value = 45264

$ws1.Range("B17").NumberFormat = "mm-dd-yy"
$ws1.Range("B17").Value = 45264

$ws1.Range("B18").NumberFormat = "mm-dd-yy"
$ws1.Range("B18").Value = 45264

# Move the active tab/selection from "Presets" to "Reference Lists".
$ws1.Activate() | Out-Null
$ws1.Range("B3").Select() | Out-Null
